# EverStation software architecture deck - "comb the DSMCC UNM DDM calling
# relationship": add a DDB box + BIOP box (with their connector arrows) below
# the existing UNM/DSI-DII/DDM chain on the DSMCC slide, nudge two existing
# connector arrows so they still line up with the UNM/DDM boxes, and bump the
# auto date placeholders on the master/layouts from 2018/10/5 to 2018/10/7.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: PowerPoint's Shape.Left/Top/Width/Height are single-precision
# points. EMU -> point -> (float32) -> EMU is lossy (floor, not round), so a
# naive conversion can land 1 EMU short of the value recorded in the OOXML.
# Search for a point value whose float32 rounding still floors back to the
# exact target EMU.
# ---------------------------------------------------------------------------
function EmuToPt([double]$emu) {
    $emuPerPt = 12700.0
    $base = $emu / $emuPerPt
    for ($i = 0; $i -lt 5000; $i++) {
        $candidate = $base + ($i * 0.0000001)
        $f32 = [float]$candidate
        $backEmu = [math]::Floor([double]$f32 * $emuPerPt)
        if ($backEmu -eq $emu) {
            return $candidate
        }
    }
    return $base
}

function SetShapeRectEmu($shape, $xEmu, $yEmu, $cxEmu, $cyEmu) {
    $shape.Left = EmuToPt $xEmu
    $shape.Top = EmuToPt $yEmu
    $shape.Width = EmuToPt $cxEmu
    $shape.Height = EmuToPt $cyEmu
}

function FindShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Slide 4 = "DSMCC的调用关系"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Nudge the two existing connector arrows (UNM -> DSI/DII, UNM -> DDM) a bit
# to the right so they keep pointing at the UNM box after the new rows push
# things around.
$conn29 = FindShapeById $s4.Shapes 29
SetShapeRectEmu $conn29 3696251 3864481 0 494367

$conn25 = FindShapeById $s4.Shapes 25
SetShapeRectEmu $conn25 3704721 4853982 0 494367

# ---------------------------------------------------------------------------
# PowerPoint hands out the lowest unused shape Id on the slide. The IDs we
# need for the five new shapes (18-22) are currently "used up" by shapes that
# were deleted earlier in the deck's history, so burn through the lower gaps
# (7,8,9,10,12,14,15) first with scratch shapes, then delete them - the Id
# counter keeps advancing instead of re-using them.
# ---------------------------------------------------------------------------
$scratch = @()
for ($i = 0; $i -lt 7; $i++) {
    $scratch += $s4.Shapes.AddTextbox(1, 0, 0, 10, 10)
}
foreach ($sh in $scratch) {
    $sh.Delete()
}

# Template shapes to clone formatting from.
$connTemplate = FindShapeById $s4.Shapes 17
$rectTemplate = FindShapeById $s4.Shapes 24

# Id 18: connector from DDM row down to the new DDB box.
$conn18 = $connTemplate.Duplicate().Item(1)
$conn18.Name = "直接连接符 17"
SetShapeRectEmu $conn18 5079818 4853982 0 494367

# Id 19: new "DDB" box.
$rect19 = $rectTemplate.Duplicate().Item(1)
$rect19.Name = "矩形 18"
SetShapeRectEmu $rect19 4408563 5374134 1426809 450572
$rect19.TextFrame.TextRange.Text = "DDB"

# Id 20: new "BIOP" box (wider, sits under UNM/DDM).
$rect20 = $rectTemplate.Duplicate().Item(1)
$rect20.Name = "矩形 19"
SetShapeRectEmu $rect20 2916584 6322400 2918788 450572
$rect20.TextFrame.TextRange.Text = "BIOP"

# Id 21: connector from DDB box down to the BIOP box.
$conn21 = $connTemplate.Duplicate().Item(1)
$conn21.Name = "直接连接符 20"
SetShapeRectEmu $conn21 5079818 5824706 0 494367

# Id 22: connector from DSI/DII box down to the BIOP box.
$conn22 = $connTemplate.Duplicate().Item(1)
$conn22.Name = "直接连接符 21"
SetShapeRectEmu $conn22 3708218 5824706 0 494367

# ---------------------------------------------------------------------------
# Bump the "last displayed" date field from 2018/10/5 to 2018/10/7 everywhere
# it appears (slide master + every slide layout).
# ---------------------------------------------------------------------------
function UpdateDateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq "2018/10/5") {
                    $sh.TextFrame.TextRange.Text = "2018/10/7"
                }
            }
        }
    }
}

UpdateDateShapes $p.SlideMaster.Shapes
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    UpdateDateShapes $layout.Shapes
}
